$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Excels ColumnWidth setter applies an internal ~5/6 character padding offset
# relative to the raw OOXML <col width> value, so compensate to land on the exact target.
$ws.Columns.Item(4).ColumnWidth = 19 - (5/6)
$ws.Columns.Item(13).ColumnWidth = 32 - (5/6)

# --- Row structure edits ---
# Insert a new row for BS170 (Q1) ahead of the resistor block.
$ws.Rows.Item(14).Insert()

# The old potentiometer row (PTV09A-2025F-B103) is being replaced entirely
# by a new part (PTD901-1015K-B103); after the insert above it now sits at row 22.
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Insert()

# Row 2 (was row 2)
$ws.Range("D2").Value = 'Arduino_Nano_v3.x'
$ws.Range("K2").Value = 1451

# Row 3 (was row 3)
$ws.Range("D3").Value = '100nF'
$ws.Range("K3").Value = 2639

# Row 4 (was row 4)
$ws.Range("A4").Value = '35ML100MEFC8X7.5'
$ws.Range("C4").Value = '1189-4176-ND'
$ws.Range("D4").Value = '100uF'
$ws.Range("I4").Value = 0.42
$ws.Range("J4").Value = '$1.26'
$ws.Range("K4").Value = 1772
$ws.Range("M4").Value = 'CAP ALUM 100UF 20% 35V RADIAL'

# Row 5 (was row 5)
$ws.Range("D5").Value = '10nF'
$ws.Range("K5").Value = 5926

# Row 6 (was row 6)
$ws.Range("D6").Value = '100nF'
$ws.Range("K6").Value = 115507

# Row 7 (was row 7)
$ws.Range("D7").Value = '1N4148'
$ws.Range("K7").Value = 245583

# Row 8 (was row 8)
$ws.Range("A8").Value = '151031VS06000'
$ws.Range("C8").Value = '732-5008-ND'
$ws.Range("D8").Value = 'LED'
$ws.Range("I8").Value = 0.15
$ws.Range("J8").Value = '$0.15'
$ws.Range("K8").Value = 16916
$ws.Range("M8").Value = 'LED GREEN DIFFUSED 3MM ROUND T/H'

# Row 9 (was row 9)
$ws.Range("D9").Value = 'AudioJack2'
$ws.Range("K9").Value = 112149

# Row 10 (was row 10)
$ws.Range("D10").Value = '2X8 Header'
$ws.Range("K10").Value = 3214

# Row 11 (was row 11)
$ws.Range("D11").Value = 'Conn_02x05_Odd_Even'
$ws.Range("K11").Value = 26153
$ws.Range("L11").Value = '9 Weeks'

# Row 12 (was row 12)
$ws.Range("D12").Value = 'Conn_01x02_Male'
$ws.Range("K12").Value = 20221

# Row 13 (was row 13)
$ws.Range("D13").Value = 'LightPipe'
$ws.Range("K13").Value = 15845

# Row 14 (new)
$ws.Range("A14").Value = 'BS170'
$ws.Range("B14").Value = 'ON Semiconductor'
$ws.Range("C14").Value = 'BS170-ND'
$ws.Range("D14").Value = 'BS170'
$ws.Range("E14").Value = 'Q1'
$ws.Range("F14").Value = 'Bulk'
$ws.Range("G14").Value = 'Active'
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0.49
$ws.Range("J14").Value = '$0.49'
$ws.Range("K14").Value = 16133
$ws.Range("L14").Value = '6 Weeks'
$ws.Range("M14").Value = 'MOSFET N-CH 60V 500MA TO-92'
$ws.Range("N14").Value = 'RoHS Compliant'
$ws.Range("O14").Value = 'Lead free'
$ws.Range("P14").Value = 'REACH Unaffected'

# Row 15 (was row 14)
$ws.Range("D15").Value = 470
$ws.Range("K15").Value = 344553

# Row 16 (was row 15)
$ws.Range("D16").Value = '10.0K'
$ws.Range("K16").Value = 5562

# Row 17 (was row 16)
$ws.Range("D17").Value = '17.8K'
$ws.Range("K17").Value = 14789

# Row 18 (was row 17)
$ws.Range("D18").Value = '140K'
$ws.Range("K18").Value = 4520

# Row 19 (was row 18)
$ws.Range("D19").Value = '15.8K'
$ws.Range("K19").Value = 3663

# Row 20 (was row 19)
$ws.Range("D20").Value = '1K'

# Row 21 (was row 20)
$ws.Range("D21").Value = '10K'
$ws.Range("K21").Value = 1273177

# Row 22 (new)
$ws.Range("A22").Value = 'PTD901-1015K-B103'
$ws.Range("B22").Value = 'Bourns Inc.'
$ws.Range("C22").Value = 'PTD901-1015K-B103-ND'
$ws.Range("D22").Value = '10K'
$ws.Range("E22").Value = 'RV1 RV2 RV3'
$ws.Range("F22").Value = 'Tray'
$ws.Range("G22").Value = 'Active'
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 1.68
$ws.Range("J22").Value = '$5.04'
$ws.Range("K22").Value = 1139
$ws.Range("L22").Value = '16 Weeks'
$ws.Range("M22").Value = 'POT 10K OHM 1/20W CARBON LINEAR'
$ws.Range("N22").Value = 'RoHS Compliant'
$ws.Range("O22").Value = 'Lead free'
$ws.Range("P22").Value = 'REACH Unaffected'

# Row 23 (was row 22)
$ws.Range("D23").Value = 'TestPoint'
$ws.Range("K23").Value = 209320

# Row 24 (was row 23)
$ws.Range("D24").Value = 'NE5532'
$ws.Range("K24").Value = 3399

# Row 25 (new)
$ws.Range("A25").Value = '1221-L'
$ws.Range("B25").Value = 'Davies Molding, LLC'
$ws.Range("C25").Value = '1722-1314-ND'
$ws.Range("D25").Value = 'Knob'
$ws.Range("E25").Value = ''
$ws.Range("F25").Value = 'Bulk'
$ws.Range("G25").Value = 'Active'
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 1.22
$ws.Range("J25").Value = '$3.66'
$ws.Range("K25").Value = 676
$ws.Range("L25").Value = '6 Weeks'
$ws.Range("M25").Value = 'KNOB SERRATED 0.236" PLASTIC'
$ws.Range("N25").Value = 'RoHS Compliant'
$ws.Range("O25").Value = 'Lead free'
$ws.Range("P25").Value = 'Not Available'

